$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.306.79'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.871.18'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '235.28'
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '0.4671'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '0.2841'
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("D9").Value = '0.06561'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '21.23'
$ws.Range("E10").Value = '  +7.58%  '
$ws.Range("D11").Value = '0.07879'
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("D12").Value = '98.21'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '1.871.10'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").Value = '5.119'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '0.6757'
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").Value = '281.22'
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").Value = '30.299.38'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '5.505'
$ws.Range("D20").Value = '12.69'
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("D21").Value = '2.113.90'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").Value = '0.000007293'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = '6.169'
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -1.19%  '
$ws.Range("D26").Value = '164.86'
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").Value = '19.18'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").Value = '1.931'
$ws.Range("E28").Value = '  -2.88%  '
$ws.Range("D29").Value = '1.374'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '0.09722'
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").Value = '4.425'
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").Value = '1.477'
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("D33").Value = '4.115'
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("D34").Value = '0.04692'
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = '1.127'
$ws.Range("E35").Value = '  +3.02%  '
$ws.Range("D36").Value = '0.7061'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '0.01862'
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = '6.292'
$ws.Range("E39").Value = '  -5.45%  '
$ws.Range("D40").Value = '2.542'
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("D41").Value = '73.41'
$ws.Range("E41").Value = '  +1.89%  '
$ws.Range("D42").Value = '1.955'
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("D43").Value = '0.8463'
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("D44").Value = '0.4179'
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '104.03'
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '7.198'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").Value = '9.136'
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("D49").Value = '933.09'
$ws.Range("E49").Value = '  -5.37%  '
$ws.Range("D50").Value = '34.09'
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("E51").Value = '  -3.08%  '
